$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 204.2
$ws.Range("I5").Value = 204.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 204.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -89.19999999999999

$ws.Range("H12").Value = 136.625
$ws.Range("I12").Value = 136.625
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 136.625
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 33.375

$ws.Range("H18").Value = 8379.583000000001
$ws.Range("I18").Value = 793.5714
$ws.Range("J18").Value = 19000
$ws.Range("K18").Value = 793.5714
$ws.Range("L18").Value = 19000
$ws.Range("M18").Value = -509.5714
$ws.Range("N18").Value = -19568

$ws.Range("H33").Value = 293.73334
$ws.Range("I33").Value = 223.53847
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 223.53847
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = 5.46153000000001
$ws.Range("N33").Value = -1208

$ws.Range("H53").Value = 844.35297
$ws.Range("I53").Value = 585.6667
$ws.Range("J53").Value = 1135.375
$ws.Range("K53").Value = 585.6667
$ws.Range("L53").Value = 1135.375
$ws.Range("M53").Value = 51.33330000000001
$ws.Range("N53").Value = -2409.375

$ws.Range("H58").Value = 658.75
$ws.Range("I58").Value = 658.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1976.25
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1826.25
$ws.Range("N58").ClearContents()

$ws.Range("H70").Value = 13593.134
$ws.Range("I70").Value = 2260
$ws.Range("J70").Value = 19259.7
$ws.Range("K70").Value = 6780
$ws.Range("L70").Value = 57779.10000000001
$ws.Range("M70").Value = -6510
$ws.Range("N70").Value = -58319.10000000001

$ws.Range("H73").Value = 13593.134
$ws.Range("I73").Value = 2260
$ws.Range("J73").Value = 19259.7
$ws.Range("K73").Value = 6780
$ws.Range("L73").Value = 57779.10000000001
$ws.Range("M73").Value = -5844
$ws.Range("N73").Value = -59651.10000000001

$ws.Range("H82").Value = 312
$ws.Range("I82").Value = 312
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 936
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -530

$ws.Range("H85").Value = 312
$ws.Range("I85").Value = 312
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 936
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 468

$ws.Range("H87").Value = 79999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 79999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495

$ws.Range("H90").Value = 79999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 79999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477

$ws.Range("H113").Value = 7577.4165
$ws.Range("I113").Value = 8829.833000000001
$ws.Range("J113").Value = 6325
$ws.Range("K113").Value = 8829.833000000001
$ws.Range("L113").Value = 6325
$ws.Range("M113").Value = -5575.833000000001
$ws.Range("N113").Value = -12833

$ws.Range("H121").Value = 4833
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 4833
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 14499
$ws.Range("N121").Value = -17993

$ws.Range("H125").Value = 5914.3335
$ws.Range("I125").Value = 5554.5557
$ws.Range("J125").Value = 6993.6665
$ws.Range("K125").Value = 49991.0013
$ws.Range("L125").Value = 62942.9985
$ws.Range("M125").Value = -47531.0013
$ws.Range("N125").Value = -67862.9985

$ws.Range("H137").Value = 35970.883
$ws.Range("I137").Value = 81518.5
$ws.Range("J137").Value = 4087.55
$ws.Range("K137").Value = 244555.5
$ws.Range("L137").Value = 12262.65
$ws.Range("M137").Value = -242005.5
$ws.Range("N137").Value = -17362.65

$ws.Range("H138").Value = 2288.7273
$ws.Range("I138").Value = 1647
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 4941
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 199
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 362.375
$ws.Range("I4").Value = 334.85715
$ws.Range("J4").Value = 555
$ws.Range("K4").Value = 334.85715
$ws.Range("L4").Value = 555
$ws.Range("M4").Value = -218.85715
$ws.Range("N4").Value = -787

$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 250
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -138
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value = 5349.8
$ws.Range("I32").Value = 2253.1843
$ws.Range("J32").Value = 31499
$ws.Range("K32").Value = 2253.1843
$ws.Range("L32").Value = 31499
$ws.Range("M32").Value = -1966.1843
$ws.Range("N32").Value = -32073

$ws.Range("H45").Value = 4999
$ws.Range("I45").Value = 4999
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4999
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4622
$ws.Range("N45").ClearContents()

$ws.Range("H102").Value = 734.1177
$ws.Range("I102").Value = 665.4545000000001
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 665.4545000000001
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 956.5454999999999
$ws.Range("N102").Value = -6244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 250
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -135
$ws.Range("N4").ClearContents()

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H22").Value = 795.1
$ws.Range("I22").Value = 800.1111
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 800.1111
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -627.1111
$ws.Range("N22").Value = -1096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 145
$ws.Range("I7").Value = 192.5
$ws.Range("J7").Value = 73.75
$ws.Range("K7").Value = 192.5
$ws.Range("L7").Value = 73.75
$ws.Range("M7").Value = -79.5
$ws.Range("N7").Value = -299.75

$ws.Range("H22").Value = 786.5
$ws.Range("I22").Value = 785.6667
$ws.Range("J22").Value = 789
$ws.Range("K22").Value = 785.6667
$ws.Range("L22").Value = 789
$ws.Range("M22").Value = -435.6667
$ws.Range("N22").Value = -1489

$ws.Range("H31").Value = 136036.6
$ws.Range("I31").Value = 202090.73
$ws.Range("J31").Value = 3928.28
$ws.Range("K31").Value = 202090.73
$ws.Range("L31").Value = 3928.28
$ws.Range("M31").Value = -201795.73
$ws.Range("N31").Value = -4518.280000000001

$ws.Range("H34").Value = 136036.6
$ws.Range("I34").Value = 202090.73
$ws.Range("J34").Value = 3928.28
$ws.Range("K34").Value = 202090.73
$ws.Range("L34").Value = 3928.28
$ws.Range("M34").Value = -201888.73
$ws.Range("N34").Value = -4332.280000000001

$ws.Range("H62").Value = 4192.8
$ws.Range("I62").Value = 4137.5713
$ws.Range("J62").Value = 4321.6665
$ws.Range("K62").Value = 4137.5713
$ws.Range("L62").Value = 4321.6665
$ws.Range("M62").Value = -3513.5713
$ws.Range("N62").Value = -5569.6665

$ws.Range("H65").Value = 4192.8
$ws.Range("I65").Value = 4137.5713
$ws.Range("J65").Value = 4321.6665
$ws.Range("K65").Value = 20687.8565
$ws.Range("L65").Value = 21608.3325
$ws.Range("M65").Value = -17567.8565
$ws.Range("N65").Value = -27848.3325

$ws.Range("H86").Value = 459135.88
$ws.Range("I86").Value = 718317.8
$ws.Range("J86").Value = 5567.5
$ws.Range("K86").Value = 718317.8
$ws.Range("L86").Value = 5567.5
$ws.Range("M86").Value = -717194.8
$ws.Range("N86").Value = -7813.5

$ws.Range("H89").Value = 459135.88
$ws.Range("I89").Value = 718317.8
$ws.Range("J89").Value = 5567.5
$ws.Range("K89").Value = 3591589
$ws.Range("L89").Value = 27837.5
$ws.Range("M89").Value = -3585973
$ws.Range("N89").Value = -39069.5

$ws.Range("H134").Value = 2249.3242
$ws.Range("I134").Value = 2314.4
$ws.Range("J134").Value = 1110.5
$ws.Range("K134").Value = 6943.200000000001
$ws.Range("L134").Value = 3331.5
$ws.Range("M134").Value = -4408.200000000001
$ws.Range("N134").Value = -8401.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4997.1113
$ws.Range("I55").Value = 1649.75
$ws.Range("J55").Value = 7675
$ws.Range("K55").Value = 4949.25
$ws.Range("L55").Value = 23025
$ws.Range("M55").Value = -4772.25
$ws.Range("N55").Value = -23379

$ws.Range("H68").Value = 2980303.8
$ws.Range("I68").Value = 10418960
$ws.Range("J68").Value = 4841.3
$ws.Range("K68").Value = 31256880
$ws.Range("L68").Value = 14523.9
$ws.Range("M68").Value = -31256069
$ws.Range("N68").Value = -16145.9

$ws.Range("H71").Value = 2980303.8
$ws.Range("I71").Value = 10418960
$ws.Range("J71").Value = 4841.3
$ws.Range("K71").Value = 93770640
$ws.Range("L71").Value = 43571.7
$ws.Range("M71").Value = -93766584
$ws.Range("N71").Value = -51683.7

$ws.Range("H113").Value = 746.4737
$ws.Range("I113").Value = 406.55554
$ws.Range("J113").Value = 1052.4
$ws.Range("K113").Value = 1219.66662
$ws.Range("L113").Value = 3157.2
$ws.Range("M113").Value = 950.33338
$ws.Range("N113").Value = -7497.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.76471
$ws.Range("I2").Value = 152.46153
$ws.Range("J2").Value = 374.75
$ws.Range("K2").Value = 152.46153
$ws.Range("L2").Value = 374.75
$ws.Range("M2").Value = -39.46153000000001
$ws.Range("N2").Value = -600.75

$ws.Range("H97").Value = 2158.652
$ws.Range("I97").Value = 2029.5
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 2029.5
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -1533.5
$ws.Range("N97").Value = -5992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 929.6667
$ws.Range("I22").Value = 929.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 929.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -634.6667
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 929.6667
$ws.Range("I27").Value = 929.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 929.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -822.6667
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2505.1365
$ws.Range("I126").Value = 2255.75
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 6767.25
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -4297.25
$ws.Range("N126").Value = -19937
